$wb = $excel.ActiveWorkbook

# Rename worksheets (sheet tab names)
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911847146533"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911873156693"
$wb.Worksheets.Item(3).Name = "RS_TO-1650291187317667"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911873636718"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911874256742"

# Sheet 1 - GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911846796513.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911846966484.csv"
$ws1.Range("B4").Value = "go_stims-16502911846986635.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911847136524.csv"

# Sheet 2 - NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16502911866206748.csv"
$ws2.Range("B3").Value = "ZB-match_6-16502911855806491.csv"
$ws2.Range("B4").Value = "OB-16502911856986487.csv"
$ws2.Range("B5").Value = "ZB-match_3-16502911852736638.csv"
$ws2.Range("B6").Value = "TB-16502911873026674.csv"
$ws2.Range("B7").Value = "TB-1650291186734665.csv"
$ws2.Range("B8").Value = "ZB-match_8-16502911850376484.csv"
$ws2.Range("B9").Value = "OB-16502911856356516.csv"
$ws2.Range("B10").Value = "TB-16502911868376791.csv"

# Sheet 4 - TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911873306684.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911873196678.csv"
$ws4.Range("B4").Value = "MM_stims-16502911873466687.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911873316674.csv"
$ws4.Range("B6").Value = "MM_stims-1650291187362672.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911873476713.csv"

# Sheet 5 - vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502911874106686.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911873956687.csv"
$ws5.Range("B4").Value = "SAT_stims-16502911873787093.csv"
$ws5.Range("B5").Value = "SAT_stims-16502911873676696.csv"
